$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1759
$ws.Range("J17").Value = 1759
$ws.Range("L17").Value = 5277
$ws.Range("N17").Value = -5613
$ws.Range("H19").Value = 1226.0526
$ws.Range("J19").Value = 1605
$ws.Range("L19").Value = 1605
$ws.Range("N19").Value = -1955
$ws.Range("H40").Value = 3627.842
$ws.Range("I40").Value = 3041.818
$ws.Range("K40").Value = 3041.818
$ws.Range("M40").Value = -2866.818
$ws.Range("H43").Value = 4822.75
$ws.Range("I43").Value = 3921.2
$ws.Range("J43").Value = 5724.3
$ws.Range("K43").Value = 3921.2
$ws.Range("L43").Value = 5724.3
$ws.Range("M43").Value = -3852.2
$ws.Range("N43").Value = -5862.3
$ws.Range("H70").Value = 2271.1428
$ws.Range("I70").Value = 2271.1428
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 6813.428400000001
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -6543.428400000001
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2271.1428
$ws.Range("I73").Value = 2271.1428
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 6813.428400000001
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -5877.428400000001
$ws.Range("N73").ClearContents()
$ws.Range("H100").Value = 1196.8889
$ws.Range("H113").Value = 83335080
$ws.Range("I113").Value = 25001872
$ws.Range("K113").Value = 25001872
$ws.Range("M113").Value = -24998618
$ws.Range("H127").Value = 9419.799999999999
$ws.Range("I127").Value = 1049.5
$ws.Range("K127").Value = 3148.5
$ws.Range("M127").Value = 1811.5
$ws.Range("H129").Value = 1232.5454
$ws.Range("J129").Value = 2749
$ws.Range("L129").Value = 8247
$ws.Range("N129").Value = -18247
$ws.Range("H135").Value = 4258.75
$ws.Range("I135").Value = 3999.5
$ws.Range("J135").Value = 4518
$ws.Range("K135").Value = 35995.5
$ws.Range("L135").Value = 40662
$ws.Range("M135").Value = -33460.5
$ws.Range("N135").Value = -45732
$ws.Range("H137").Value = 3796.1836
$ws.Range("I137").Value = 2756.3333
$ws.Range("J137").Value = 5940.875
$ws.Range("K137").Value = 8268.999899999999
$ws.Range("L137").Value = 17822.625
$ws.Range("M137").Value = -5718.999899999999
$ws.Range("N137").Value = -22922.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 451.375
$ws.Range("I5").Value = 176.5
$ws.Range("J5").Value = 726.25
$ws.Range("K5").Value = 176.5
$ws.Range("L5").Value = 726.25
$ws.Range("M5").Value = -64.5
$ws.Range("N5").Value = -950.25
$ws.Range("H32").Value = 8337754.5
$ws.Range("I32").Value = 10206030
$ws.Range("K32").Value = 10206030
$ws.Range("M32").Value = -10205743
$ws.Range("H88").Value = 4030.3215
$ws.Range("I88").Value = 3170.8333
$ws.Range("K88").Value = 3170.8333
$ws.Range("M88").Value = -2764.8333
$ws.Range("H91").Value = 4030.3215
$ws.Range("I91").Value = 3170.8333
$ws.Range("K91").Value = 3170.8333
$ws.Range("M91").Value = -1766.8333
$ws.Range("H132").Value = 4114.879
$ws.Range("I132").Value = 1796.4762
$ws.Range("K132").Value = 5389.4286
$ws.Range("M132").Value = -2859.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 451.375
$ws.Range("I4").Value = 176.5
$ws.Range("J4").Value = 726.25
$ws.Range("K4").Value = 176.5
$ws.Range("L4").Value = 726.25
$ws.Range("M4").Value = -61.5
$ws.Range("N4").Value = -956.25
$ws.Range("H105").Value = 3133.25
$ws.Range("I105").Value = 5700
$ws.Range("K105").Value = 5700
$ws.Range("M105").Value = -3953
$ws.Range("H107").Value = 3258.25
$ws.Range("I107").Value = 3258.25
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3258.25
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1338.25
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 1114415.4
$ws.Range("I134").Value = 3435
$ws.Range("K134").Value = 10305
$ws.Range("M134").Value = -7770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 873.5294
$ws.Range("I16").Value = 672.8
$ws.Range("J16").Value = 1160.2858
$ws.Range("K16").Value = 672.8
$ws.Range("L16").Value = 1160.2858
$ws.Range("M16").Value = -385.8
$ws.Range("N16").Value = -1734.2858
$ws.Range("H99").Value = 3429.9375
$ws.Range("I99").Value = 3123.25
$ws.Range("J99").Value = 3736.625
$ws.Range("K99").Value = 3123.25
$ws.Range("L99").Value = 3736.625
$ws.Range("M99").Value = -1625.25
$ws.Range("N99").Value = -6732.625
$ws.Range("H113").Value = 873.5294
$ws.Range("I113").Value = 672.8
$ws.Range("J113").Value = 1160.2858
$ws.Range("K113").Value = 672.8
$ws.Range("L113").Value = 1160.2858
$ws.Range("M113").Value = 1497.2
$ws.Range("N113").Value = -5500.2858
$ws.Range("H126").Value = 3429.9375
$ws.Range("I126").Value = 3123.25
$ws.Range("J126").Value = 3736.625
$ws.Range("K126").Value = 9369.75
$ws.Range("L126").Value = 11209.875
$ws.Range("M126").Value = -6899.75
$ws.Range("N126").Value = -16149.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2014.6428
$ws.Range("I12").Value = 3746.3333
$ws.Range("J12").Value = 715.875
$ws.Range("K12").Value = 11238.9999
$ws.Range("L12").Value = 2147.625
$ws.Range("M12").Value = -11065.9999
$ws.Range("N12").Value = -2493.625
$ws.Range("H109").Value = 2825.6667
$ws.Range("J109").Value = 2700
$ws.Range("L109").Value = 8100
$ws.Range("N109").Value = -10180
$ws.Range("H131").Value = 7434.875
$ws.Range("J131").Value = 7573.1914
$ws.Range("L131").Value = 22719.5742
$ws.Range("N131").Value = -32799.5742
$ws.Range("H139").Value = 3316.3333
$ws.Range("I139").Value = 1891.9231
$ws.Range("K139").Value = 5675.7693
$ws.Range("M139").Value = -535.7692999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 2050101.8
$ws.Range("I7").Value = 5000250
$ws.Range("J7").Value = 83336.336
$ws.Range("K7").Value = 5000250
$ws.Range("L7").Value = 83336.336
$ws.Range("M7").Value = -5000138
$ws.Range("N7").Value = -83560.336
$ws.Range("H8").Value = 2050101.8
$ws.Range("I8").Value = 5000250
$ws.Range("J8").Value = 83336.336
$ws.Range("K8").Value = 5000250
$ws.Range("L8").Value = 83336.336
$ws.Range("M8").Value = -5000111
$ws.Range("N8").Value = -83614.336
$ws.Range("H70").Value = 21667.666
$ws.Range("I70").Value = 24001.2
$ws.Range("K70").Value = 24001.2
$ws.Range("M70").Value = -23731.2
$ws.Range("H73").Value = 21667.666
$ws.Range("I73").Value = 24001.2
$ws.Range("K73").Value = 24001.2
$ws.Range("M73").Value = -23065.2
$ws.Range("H97").Value = 1785.8636
$ws.Range("I97").Value = 1810.4667
$ws.Range("J97").Value = 1733.1428
$ws.Range("K97").Value = 1810.4667
$ws.Range("L97").Value = 1733.1428
$ws.Range("M97").Value = -1314.4667
$ws.Range("N97").Value = -2725.1428
$ws.Range("H126").Value = 4655.636
$ws.Range("I126").Value = 4303
$ws.Range("J126").Value = 4857.143
$ws.Range("K126").Value = 12909
$ws.Range("L126").Value = 14571.429
$ws.Range("M126").Value = -10439
$ws.Range("N126").Value = -19511.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1553.0667
$ws.Range("I22").Value = 1526.7273
$ws.Range("K22").Value = 1526.7273
$ws.Range("M22").Value = -1231.7273
$ws.Range("H27").Value = 1553.0667
$ws.Range("I27").Value = 1526.7273
$ws.Range("K27").Value = 1526.7273
$ws.Range("M27").Value = -1419.7273
$ws.Range("H55").Value = 47619696
$ws.Range("J55").Value = 582.5
$ws.Range("L55").Value = 582.5
$ws.Range("N55").Value = -928.5
$ws.Range("H132").Value = 567256.6
$ws.Range("I132").Value = 14186.714
$ws.Range("J132").Value = 2503001.2
$ws.Range("K132").Value = 42560.142
$ws.Range("L132").Value = 7509003.600000001
$ws.Range("M132").Value = -40030.142
$ws.Range("N132").Value = -7514063.600000001
$ws.Range("H136").Value = 78040.21000000001
$ws.Range("I136").Value = 12975.4
$ws.Range("K136").Value = 38926.2
$ws.Range("M136").Value = -36376.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 46660
$ws.Range("J80").Value = 80000
$ws.Range("L80").Value = 80000
$ws.Range("N80").Value = -81996
$ws.Range("H83").Value = 46660
$ws.Range("J83").Value = 80000
$ws.Range("L83").Value = 240000
$ws.Range("N83").Value = -249984
$ws.Range("H123").Value = 75030
$ws.Range("J123").Value = 75030
$ws.Range("L123").Value = 75030
$ws.Range("N123").Value = -84830
